$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (D) and volume-change (E) columns to refreshed values.
# D-column values are forced to text (NumberFormat "@") before assignment and then
# ClearFormats() is used to drop the temporary style again, so the cells keep their
# original un-styled look while the literal text (e.g. trailing zeros) is preserved.

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "26.780.25"
$cell.ClearFormats()
$ws.Range("E2").Value = "  +7.81%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.742.82"
$cell.ClearFormats()
$ws.Range("E3").Value = "  +4.53%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.ClearFormats()
$ws.Range("E4").Value = "  -0.43%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "335.28"
$cell.ClearFormats()
$ws.Range("E5").Value = "  +1.93%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.9995"
$cell.ClearFormats()
$ws.Range("E6").Value = "  -0.29%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.3751"
$cell.ClearFormats()
$ws.Range("E7").Value = "  +2.96%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "48.47"
$cell.ClearFormats()
$ws.Range("E8").Value = "  +4.51%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.3387"
$cell.ClearFormats()
$ws.Range("E9").Value = "  +5.04%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "1.191"
$cell.ClearFormats()
$ws.Range("E10").Value = "  +4.88%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.07477"
$cell.ClearFormats()
$ws.Range("E11").Value = "  +6.43%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "1.000"
$cell.ClearFormats()
$ws.Range("E12").Value = "  -0.19%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "6.401"
$cell.ClearFormats()
$ws.Range("E13").Value = "  +6.00%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "20.41"
$cell.ClearFormats()
$ws.Range("E14").Value = "  +5.02%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "7.061"
$cell.ClearFormats()
$ws.Range("E15").Value = "  +7.24%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "1.741.38"
$cell.ClearFormats()
$ws.Range("E16").Value = "  +4.13%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.00001080"
$cell.ClearFormats()
$ws.Range("E17").Value = "  +3.62%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.06730"
$cell.ClearFormats()
$ws.Range("E18").Value = "  +2.56%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "82.79"
$cell.ClearFormats()
$ws.Range("E19").Value = "  +5.62%  "

$ws.Range("E20").Value = "  -0.12%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "16.73"
$cell.ClearFormats()
$ws.Range("E21").Value = "  +6.24%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "6.233"
$cell.ClearFormats()
$ws.Range("E22").Value = "  +5.83%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "12.78"
$cell.ClearFormats()
$ws.Range("E23").Value = "  -0.52%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "26.759.01"
$cell.ClearFormats()
$ws.Range("E24").Value = "  +7.60%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.465"
$cell.ClearFormats()
$ws.Range("E25").Value = "  +0.97%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "1.475"
$cell.ClearFormats()
$ws.Range("E26").Value = "  +26.14%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "2.415"
$cell.ClearFormats()
$ws.Range("E27").Value = "  +1.84%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "152.76"
$cell.ClearFormats()
$ws.Range("E28").Value = "  +3.43%  "

$ws.Range("E29").Value = "  +5.49%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "1.936.23"
$cell.ClearFormats()
$ws.Range("E30").Value = "  +3.97%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "132.38"
$cell.ClearFormats()
$ws.Range("E31").Value = "  +5.63%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "4.120"
$cell.ClearFormats()
$ws.Range("E32").Value = "  +1.06%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "6.033"
$cell.ClearFormats()
$ws.Range("E33").Value = "  +5.75%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.08627"
$cell.ClearFormats()
$ws.Range("E34").Value = "  +2.36%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.688"
$cell.ClearFormats()
$ws.Range("E35").Value = "  +2.63%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "12.91"
$cell.ClearFormats()
$ws.Range("E36").Value = "  +6.25%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "5.432"
$cell.ClearFormats()
$ws.Range("E37").Value = "  +5.85%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.02357"
$cell.ClearFormats()
$ws.Range("E38").Value = "  +5.56%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.2176"
$cell.ClearFormats()
$ws.Range("E39").Value = "  +4.66%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.06266"
$cell.ClearFormats()
$ws.Range("E40").Value = "  +4.77%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "8.494"
$cell.ClearFormats()
$ws.Range("E41").Value = "  +4.05%  "

$ws.Range("E42").Value = "  -0.37%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.6260"
$cell.ClearFormats()
$ws.Range("E43").Value = "  +6.01%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "14.35"
$cell.ClearFormats()
$ws.Range("E44").Value = "  +4.84%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "1.0000"
$cell.ClearFormats()
$ws.Range("E45").Value = "  -0.05%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "3.929"
$cell.ClearFormats()
$ws.Range("E46").Value = "  +2.42%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.6074"
$cell.ClearFormats()
$ws.Range("E47").Value = "  +6.47%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "129.31"
$cell.ClearFormats()
$ws.Range("E48").Value = "  +4.07%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "2.073"
$cell.ClearFormats()
$ws.Range("E49").Value = "  +6.33%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.07222"
$cell.ClearFormats()
$ws.Range("E50").Value = "  +3.46%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "77.86"
$cell.ClearFormats()
$ws.Range("E51").Value = "  +4.94%  "

